$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 202.44444
$ws.Range("I33").Value = 198.93333
$ws.Range("J33").Value = 220
$ws.Range("K33").Value = 198.93333
$ws.Range("L33").Value = 220
$ws.Range("M33").Value = 30.06666999999999
$ws.Range("N33").Value = -678
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H87").Value = 54327.5
$ws.Range("J87").Value = 54327.5
$ws.Range("L87").Value = 54327.5
$ws.Range("N87").Value = -56823.5
$ws.Range("H88").Value = 1485.5
$ws.Range("I88").Value = 1650
$ws.Range("J88").Value = 1321
$ws.Range("K88").Value = 1650
$ws.Range("L88").Value = 1321
$ws.Range("M88").Value = -1244
$ws.Range("N88").Value = -2133
$ws.Range("H90").Value = 54327.5
$ws.Range("J90").Value = 54327.5
$ws.Range("L90").Value = 162982.5
$ws.Range("N90").Value = -175462.5
$ws.Range("H91").Value = 1485.5
$ws.Range("I91").Value = 1650
$ws.Range("J91").Value = 1321
$ws.Range("K91").Value = 1650
$ws.Range("L91").Value = 1321
$ws.Range("M91").Value = -246
$ws.Range("N91").Value = -4129
$ws.Range("H92").Value = 737.3333
$ws.Range("J92").Value = 499
$ws.Range("L92").Value = 499
$ws.Range("N92").Value = -2995
$ws.Range("H94").Value = 1733.3334
$ws.Range("I94").Value = 1333.3334
$ws.Range("K94").Value = 1333.3334
$ws.Range("M94").Value = -882.3334
$ws.Range("H96").Value = 399.26666
$ws.Range("I96").Value = 296.16666
$ws.Range("J96").Value = 811.6667
$ws.Range("K96").Value = 888.4999799999999
$ws.Range("L96").Value = 2435.0001
$ws.Range("M96").Value = 484.5000200000001
$ws.Range("N96").Value = -5181.0001
$ws.Range("H100").Value = 3304.8
$ws.Range("J100").Value = 1699.5
$ws.Range("L100").Value = 1699.5
$ws.Range("N100").Value = -2781.5
$ws.Range("H116").Value = 5272.5
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H125").Value = 505.33334
$ws.Range("I125").Value = 406.4
$ws.Range("K125").Value = 3657.6
$ws.Range("M125").Value = -1197.6
$ws.Range("H131").Value = 847.5
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 995
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2985
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -13065
$ws.Range("H137").Value = 1028.7693
$ws.Range("I137").Value = 803.44446
$ws.Range("K137").Value = 2410.33338
$ws.Range("M137").Value = 139.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17500
$ws.Range("J24").Value = 17500
$ws.Range("L24").Value = 17500
$ws.Range("N24").Value = -18248
$ws.Range("H32").Value = 21600.133
$ws.Range("I32").Value = 21000.143
$ws.Range("K32").Value = 21000.143
$ws.Range("M32").Value = -20713.143
$ws.Range("H100").Value = 17500
$ws.Range("J100").Value = 17500
$ws.Range("L100").Value = 17500
$ws.Range("N100").Value = -19664
$ws.Range("H102").Value = 951.13336
$ws.Range("I102").Value = 951.13336
$ws.Range("K102").Value = 951.13336
$ws.Range("M102").Value = 670.86664
$ws.Range("H132").Value = 1156.9025
$ws.Range("I132").Value = 881.0833
$ws.Range("K132").Value = 2643.2499
$ws.Range("M132").Value = -113.2498999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2223.3333
$ws.Range("I86").Value = 1964.4445
$ws.Range("K86").Value = 1964.4445
$ws.Range("M86").Value = -841.4445000000001
$ws.Range("H89").Value = 2223.3333
$ws.Range("I89").Value = 1964.4445
$ws.Range("K89").Value = 9822.2225
$ws.Range("M89").Value = -4206.2225
$ws.Range("H134").Value = 2268.15
$ws.Range("I134").Value = 2124.3684
$ws.Range("K134").Value = 6373.1052
$ws.Range("M134").Value = -3838.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 31899.8
$ws.Range("J95").Value = 31899.8
$ws.Range("L95").Value = 31899.8
$ws.Range("N95").Value = -37391.8
$ws.Range("H105").Value = 4750
$ws.Range("I105").Value = 4750
$ws.Range("K105").Value = 4750
$ws.Range("M105").Value = -3003
$ws.Range("H122").Value = 3777
$ws.Range("I122").Value = 3856.1428
$ws.Range("K122").Value = 11568.4284
$ws.Range("M122").Value = -9118.428400000001
$ws.Range("H132").Value = 3397.9375
$ws.Range("I132").Value = 2420.3333
$ws.Range("K132").Value = 7260.999899999999
$ws.Range("M132").Value = -4730.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 8560
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 8400
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = -7714
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 8560
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 25200
$ws.Range("L65").Value = 90000
$ws.Range("M65").Value = -21768
$ws.Range("N65").Value = -96864
$ws.Range("H80").Value = 14000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 14000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 42000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -43872
$ws.Range("H83").Value = 14000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 14000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 126000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -135360
$ws.Range("H98").Value = 2110.6
$ws.Range("I98").Value = 1390
$ws.Range("J98").Value = 3191.5
$ws.Range("K98").Value = 4170
$ws.Range("L98").Value = 9574.5
$ws.Range("M98").Value = -2672
$ws.Range("N98").Value = -12570.5
$ws.Range("H136").Value = 1974.2
$ws.Range("I136").Value = 1842.75
$ws.Range("K136").Value = 5528.25
$ws.Range("M136").Value = -428.25
$ws.Range("H137").Value = 2429.8462
$ws.Range("I137").Value = 2105.4546
$ws.Range("J137").Value = 4214
$ws.Range("K137").Value = 6316.3638
$ws.Range("L137").Value = 12642
$ws.Range("M137").Value = -1216.3638
$ws.Range("N137").Value = -22842
$ws.Range("H138").Value = 2800
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280
$ws.Range("H139").Value = 4000
$ws.Range("I139").Value = 2000
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 6000
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -860
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3276.9285
$ws.Range("I132").Value = 2717
$ws.Range("J132").Value = 3588
$ws.Range("K132").Value = 8151
$ws.Range("L132").Value = 10764
$ws.Range("M132").Value = -5621
$ws.Range("N132").Value = -15824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2428.8
$ws.Range("I16").Value = 2925
$ws.Range("J16").Value = 444
$ws.Range("K16").Value = 2925
$ws.Range("L16").Value = 444
$ws.Range("M16").Value = -2755
$ws.Range("N16").Value = -784
$ws.Range("H22").Value = 8915.666999999999
$ws.Range("I22").Value = 9695.5
$ws.Range("K22").Value = 9695.5
$ws.Range("M22").Value = -9400.5
$ws.Range("H27").Value = 8915.666999999999
$ws.Range("I27").Value = 9695.5
$ws.Range("K27").Value = 9695.5
$ws.Range("M27").Value = -9588.5
$ws.Range("H100").Value = 2399.4443
$ws.Range("J100").Value = 2950
$ws.Range("L100").Value = 2950
$ws.Range("N100").Value = -4032
$ws.Range("H136").Value = 3788
$ws.Range("I136").Value = 3788
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11364
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8814
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 49997.5
$ws.Range("J74").Value = 49997.5
$ws.Range("L74").Value = 49997.5
$ws.Range("N74").Value = -51869.5
$ws.Range("H77").Value = 49997.5
$ws.Range("J77").Value = 49997.5
$ws.Range("L77").Value = 149992.5
$ws.Range("N77").Value = -159352.5
$ws.Range("H96").Value = 1799.5
$ws.Range("J96").Value = 1699.5
$ws.Range("L96").Value = 1699.5
$ws.Range("N96").Value = -4445.5
$ws.Range("H104").Value = 40185
$ws.Range("J104").Value = 40185
$ws.Range("L104").Value = 40185
$ws.Range("N104").Value = -47173
$ws.Range("H122").Value = 1120.625
$ws.Range("I122").Value = 1120.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3361.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -911.875
$ws.Range("N122").ClearContents()
